# Generate Report for Handoff
# Adds two new file entries (99c0cb08-0261-46ce-b817-9371000eeb7e and
# f25d11b0-9aee-4bf0-a370-8d13c5afb789) around the existing
# c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc row on every sheet:
#  - Overview sheet (row per file)
#  - zh-cn sheet (row per file)
#  - de-de sheet (row per file)
# The new rows are inserted so the file ordering stays alphabetical-ish,
# matching the original commit: 99c0cb08 goes in right before c0f34807,
# and f25d11b0 goes in right after it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Insert a row above the current c0f34807 row (row 5) for 99c0cb08,
# then a row below the (now shifted) c0f34807 row (row 7) for f25d11b0.
$ws1.Rows.Item(5).Insert()
$ws1.Rows.Item(7).Insert()

$ws1.Range("A5").Value = "99c0cb08-0261-46ce-b817-9371000eeb7e.md"
$ws1.Range("B5").Value = "e2e\99c0cb08-0261-46ce-b817-9371000eeb7e.md"
$ws1.Range("C5").Value = ".md"
$ws1.Range("E5").Value = "Ready for handoff"
$ws1.Range("F5").Value = "Ready for handoff"
$ws1.Range("G5").Value = "2016-08-26 22:40:39"

$ws1.Range("A7").Value = "f25d11b0-9aee-4bf0-a370-8d13c5afb789.md"
$ws1.Range("B7").Value = "e2e\f25d11b0-9aee-4bf0-a370-8d13c5afb789.md"
$ws1.Range("C7").Value = ".md"
$ws1.Range("E7").Value = "Ready for handoff"
$ws1.Range("F7").Value = "Ready for handoff"
$ws1.Range("G7").Value = "2016-08-26 22:40:39"

# Resize the "Overview" table to cover the two new rows.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G7"))

# Rebuild the hyperlinks on column B in document order so relationship ids
# come out as rId2..rId7 (rId1 is the table part).
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c9d2eafe0cc187b420d55dfb0d6e4caf0ec3b3d/e2e/00b79f86-3c1f-43da-b881-be20b8c858de.md", "", "", "e2e\00b79f86-3c1f-43da-b881-be20b8c858de.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0df7bcbfc69146c66704463055548ee361e545eb/e2e/6ee639aa-19b3-4a69-a0d4-c1158b77f850.md", "", "", "e2e\6ee639aa-19b3-4a69-a0d4-c1158b77f850.md")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0df7bcbfc69146c66704463055548ee361e545eb/e2e/e687fb83-a44d-4904-b82d-23191b02eef8.md", "", "", "e2e\e687fb83-a44d-4904-b82d-23191b02eef8.md")
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99c0cb08acde1234567890abcdef1234567890ab/e2e/99c0cb08-0261-46ce-b817-9371000eeb7e.md", "", "", "e2e\99c0cb08-0261-46ce-b817-9371000eeb7e.md")
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40fbf96366ecc672b4088f06e0367da3e8cdf7c5/e2e/c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md", "", "", "e2e\c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md")
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f25d11b0acde1234567890abcdef1234567890ab/e2e/f25d11b0-9aee-4bf0-a370-8d13c5afb789.md", "", "", "e2e\f25d11b0-9aee-4bf0-a370-8d13c5afb789.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(5).Insert()
$ws2.Rows.Item(7).Insert()

$ws2.Range("A5").Value = "99c0cb08-0261-46ce-b817-9371000eeb7e.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "e2e"
$ws2.Range("E5").Value = "ht"
$ws2.Range("F5").Value = "False"
$ws2.Range("G5").Value = "99c0cb08-0261-46ce-b817-9371000eeb7e.ce0acbc10a171cbaa61d56f2547d01301e2a16e2.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-26 22:40:34"
$ws2.Range("K5").Value = "0001-01-01 00:00:00"
$ws2.Range("M5").Value = "True"
$ws2.Range("O5").Value = "False"

$ws2.Range("A7").Value = "f25d11b0-9aee-4bf0-a370-8d13c5afb789.md"
$ws2.Range("B7").Value = ".md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "e2e"
$ws2.Range("E7").Value = "ht"
$ws2.Range("F7").Value = "False"
$ws2.Range("G7").Value = "f25d11b0-9aee-4bf0-a370-8d13c5afb789.765c1cfd17a176be099ea7213e16ff2655bd99b7.zh-cn.xlf"
$ws2.Range("H7").Value = "2016-08-26 22:40:34"
$ws2.Range("K7").Value = "0001-01-01 00:00:00"
$ws2.Range("M7").Value = "True"
$ws2.Range("O7").Value = "False"

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P7"))

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c9d2eafe0cc187b420d55dfb0d6e4caf0ec3b3d/e2e/00b79f86-3c1f-43da-b881-be20b8c858de.md", "", "", "00b79f86-3c1f-43da-b881-be20b8c858de.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/148d272ff01e6e7a62894d4ea1b681b3bf826e15/e2e/00b79f86-3c1f-43da-b881-be20b8c858de.md", "", "", "00b79f86-3c1f-43da-b881-be20b8c858de.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0df7bcbfc69146c66704463055548ee361e545eb/e2e/6ee639aa-19b3-4a69-a0d4-c1158b77f850.md", "", "", "6ee639aa-19b3-4a69-a0d4-c1158b77f850.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0df7bcbfc69146c66704463055548ee361e545eb/e2e/e687fb83-a44d-4904-b82d-23191b02eef8.md", "", "", "e687fb83-a44d-4904-b82d-23191b02eef8.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99c0cb08acde1234567890abcdef1234567890ab/e2e/99c0cb08-0261-46ce-b817-9371000eeb7e.md", "", "", "99c0cb08-0261-46ce-b817-9371000eeb7e.md")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40fbf96366ecc672b4088f06e0367da3e8cdf7c5/e2e/c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md", "", "", "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f25d11b0acde1234567890abcdef1234567890ab/e2e/f25d11b0-9aee-4bf0-a370-8d13c5afb789.md", "", "", "f25d11b0-9aee-4bf0-a370-8d13c5afb789.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(5).Insert()
$ws3.Rows.Item(7).Insert()

$ws3.Range("A5").Value = "99c0cb08-0261-46ce-b817-9371000eeb7e.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "e2e"
$ws3.Range("E5").Value = "ht"
$ws3.Range("F5").Value = "False"
$ws3.Range("G5").Value = "99c0cb08-0261-46ce-b817-9371000eeb7e.ce0acbc10a171cbaa61d56f2547d01301e2a16e2.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-26 22:40:39"
$ws3.Range("K5").Value = "0001-01-01 00:00:00"
$ws3.Range("M5").Value = "True"
$ws3.Range("O5").Value = "False"

$ws3.Range("A7").Value = "f25d11b0-9aee-4bf0-a370-8d13c5afb789.md"
$ws3.Range("B7").Value = ".md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "e2e"
$ws3.Range("E7").Value = "ht"
$ws3.Range("F7").Value = "False"
$ws3.Range("G7").Value = "f25d11b0-9aee-4bf0-a370-8d13c5afb789.765c1cfd17a176be099ea7213e16ff2655bd99b7.de-de.xlf"
$ws3.Range("H7").Value = "2016-08-26 22:40:39"
$ws3.Range("K7").Value = "0001-01-01 00:00:00"
$ws3.Range("M7").Value = "True"
$ws3.Range("O7").Value = "False"

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P7"))

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c9d2eafe0cc187b420d55dfb0d6e4caf0ec3b3d/e2e/00b79f86-3c1f-43da-b881-be20b8c858de.md", "", "", "00b79f86-3c1f-43da-b881-be20b8c858de.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d9af93a35ffb34a58c6d042b6c555bc9d2dae411/e2e/00b79f86-3c1f-43da-b881-be20b8c858de.md", "", "", "00b79f86-3c1f-43da-b881-be20b8c858de.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0df7bcbfc69146c66704463055548ee361e545eb/e2e/6ee639aa-19b3-4a69-a0d4-c1158b77f850.md", "", "", "6ee639aa-19b3-4a69-a0d4-c1158b77f850.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0df7bcbfc69146c66704463055548ee361e545eb/e2e/e687fb83-a44d-4904-b82d-23191b02eef8.md", "", "", "e687fb83-a44d-4904-b82d-23191b02eef8.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99c0cb08acde1234567890abcdef1234567890ab/e2e/99c0cb08-0261-46ce-b817-9371000eeb7e.md", "", "", "99c0cb08-0261-46ce-b817-9371000eeb7e.md")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40fbf96366ecc672b4088f06e0367da3e8cdf7c5/e2e/c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md", "", "", "c0f34807-8fa5-476c-8a5d-1ce4e0bbe3fc.md")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f25d11b0acde1234567890abcdef1234567890ab/e2e/f25d11b0-9aee-4bf0-a370-8d13c5afb789.md", "", "", "f25d11b0-9aee-4bf0-a370-8d13c5afb789.md")

Write-Host "Done"
